$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''42.314.11'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  +0.26%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.282.94'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.45%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '''  +0.94%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''311.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -1.42%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''101.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.54%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -1.19%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''1.00'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  +0.22%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''0.594'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.72%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''38.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  -2.53%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -0.89%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''8.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -2.26%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '''0.109'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '''  +1.92%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''0.976'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  +1.78%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''15.02'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '''  -1.01%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '''2.624.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -0.57%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''2.314.90'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  +1.69%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''42.533.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  +0.76%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''7.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -1.53%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '''  -0.92%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''13.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  +12.03%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''72.89'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.55%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''3.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '''  -0.73%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''262.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '''  -5.03%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''2.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  -3.97%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  +0.07%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''10.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.65%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = '''2.33'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  -1.41%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''6.87'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +15.14%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = '''22.29'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '''  -2.07%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''36.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -3.49%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = '''165.72'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '''  -0.12%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''0.0861'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  -1.32%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = '''WEMIXToken'
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = '''2.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '''  -1.46%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = '''Stellar'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = '''0.130'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -3.14%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.111'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '''  -4.86%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''4.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '''  -1.82%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = '''  -4.37%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''3.67'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -0.70%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''2.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -3.80%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''1.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +4.86%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''69.12'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -1.03%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = '''Algorand'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = '''0.227'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '''  +0.65%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("B44").Value = '''FirstDigitalUSD'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").Value = '''1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '''  +0.41%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '''94.06'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -2.08%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = '''Celestia'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = '''https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = '''11.94'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -0.61%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = '''Maker'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = '''1.715.64'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  +7.82%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = '''ordi'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = '''https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = '''78.89'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -0.94%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = '''Aave'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = '''https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = '''110.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '''  -2.28%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = '''FraxShare'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = '''8.67'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -3.40%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = '''THORChain'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = '''https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = '''5.16'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -1.97%  '
$ws.Range("E51").Style = "Normal"
